$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.01.2022 15:45"

# D2: change from text "+0.6" to a real number 0.6
$ws.Range("D2").Value = 0.6

# E2: change from text timestamp to a real date serial number,
# using the same date/time number format as the other rows (E3:E10)
$ws.Range("E2").Value = 44582.64587962963
$ws.Range("E2").NumberFormat = $ws.Range("E3").NumberFormat
